# Generate Report for Handback
# --------------------------------------------------------------------------
# The 7bf84cc6-... item has now been handed back (in sync with en-US), so it
# moves to the top of the report (row 2) with a new status, while the
# 05ea4a7e-... item (still "Ready for handoff") drops to row 3. The two
# per-language sheets also grow two columns ("Latest Target File" /
# "Latest Handback File") that get populated for the handed-back item.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Common hyperlink targets reused across sheets.
$md7bf8Url   = "https://github.com/OpenLocalizationTest/oltest/blob/d6399064f0bb4628bb80159e2e7cdfb1fef0a03c/e2e/7bf84cc6-557b-4a5e-9688-f355a2432ea0.md"
$md05eaUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/afb5be80dae66df280bb7df037e1aef29fad204e/e2e/05ea4a7e-cab7-444d-9b45-abe39f53d2f7.md"
$configUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/d6399064f0bb4628bb80159e2e7cdfb1fef0a03c/.localization-config"

$md7bf8Name  = "7bf84cc6-557b-4a5e-9688-f355a2432ea0.md"
$md05eaName  = "05ea4a7e-cab7-444d-9b45-abe39f53d2f7.md"
$configName  = ".localization-config"

$handedBack  = "Handed back: in sync with en-US"
$readyOff    = "Ready for handoff"
$notLocalized = "Not to be localized"

# ===========================================================================
# Sheet 1: Overview
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop every existing hyperlink on the sheet (this API clears the whole
# sheet's collection regardless of which range it is invoked on) so we can
# rebuild them cleanly in the new row order.
$wsOverview.Range("A1").Hyperlinks.Delete()

$wsOverview.Range("A2").Value = $md7bf8Name
$wsOverview.Range("B2").Value = $handedBack
$wsOverview.Range("C2").Value = $handedBack

$wsOverview.Range("A3").Value = $md05eaName
$wsOverview.Range("B3").Value = $readyOff
$wsOverview.Range("C3").Value = $readyOff

$wsOverview.Range("A4").Value = $configName
$wsOverview.Range("B4").Value = $notLocalized
$wsOverview.Range("C4").Value = $notLocalized

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $md7bf8Url, "", "", $md7bf8Name)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $md05eaUrl, "", "", $md05eaName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $configUrl, "", "", $configName)

$wsOverview.Range("A2").Style = "HyperLink"
$wsOverview.Range("A3").Style = "HyperLink"
$wsOverview.Range("A4").Style = "HyperLink"

# ===========================================================================
# Sheet 2: zh-cn
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$zh7bf8XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/75b8123256fe994335c5f73a8029c4703a0979bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7bf84cc6-557b-4a5e-9688-f355a2432ea0.c0005c834c23b3316dcc6d7d8b2418b1153ebe82.zh-cn.xlf"
$zh05eaXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/45d4e2bb1ac43240f744f7025c83a52485a64a2d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/05ea4a7e-cab7-444d-9b45-abe39f53d2f7.544af659032ec6d29b8d61d45c87655c2dbfde8c.zh-cn.xlf"
$zh7bf8XlfName = "7bf84cc6-557b-4a5e-9688-f355a2432ea0.c0005c834c23b3316dcc6d7d8b2418b1153ebe82.zh-cn.xlf"
$zh05eaXlfName = "05ea4a7e-cab7-444d-9b45-abe39f53d2f7.544af659032ec6d29b8d61d45c87655c2dbfde8c.zh-cn.xlf"

$wsZh.Range("A1").Hyperlinks.Delete()

# Row 2: 7bf84cc6 item - handed back, fully populated (incl. new columns).
$wsZh.Range("A2").Value = $md7bf8Name
$wsZh.Range("B2").Value = $handedBack
$wsZh.Range("C2").Value = $zh7bf8XlfName
$wsZh.Range("D2").Value = "2016-03-11 00:50:36"
$wsZh.Range("E2").Value = $md7bf8Name
$wsZh.Range("F2").Value = $zh7bf8XlfName
$wsZh.Range("G2").Value = "2016-03-11 00:51:28"
$wsZh.Range("H2").Value = "Include"

# Row 3: 05ea4a7e item - still ready for handoff.
$wsZh.Range("A3").Value = $md05eaName
$wsZh.Range("B3").Value = $readyOff
$wsZh.Range("C3").Value = $zh05eaXlfName
$wsZh.Range("D3").Value = "2016-03-11 00:50:06"
$wsZh.Range("G3").Value = "0001-01-01 00:00:00"
$wsZh.Range("H3").Value = "Include"

# Row 4: .localization-config - unchanged.
$wsZh.Range("A4").Value = $configName
$wsZh.Range("B4").Value = $notLocalized
$wsZh.Range("D4").Value = "0001-01-01 00:00:00"
$wsZh.Range("G4").Value = "0001-01-01 00:00:00"
$wsZh.Range("H4").Value = "Ignored"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $md7bf8Url, "", "", $md7bf8Name)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zh7bf8XlfUrl, "", "", $zh7bf8XlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $md7bf8Url, "", "", $md7bf8Name)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zh7bf8XlfUrl, "", "", $zh7bf8XlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $md05eaUrl, "", "", $md05eaName)
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), $zh05eaXlfUrl, "", "", $zh05eaXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $configUrl, "", "", $configName)

foreach ($addr in @("A2", "C2", "E2", "F2", "A3", "C3", "A4")) {
    $wsZh.Range($addr).Style = "HyperLink"
}
$wsZh.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ===========================================================================
# Sheet 3: de-de
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$de7bf8XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c32577e01d9263c9978c4522bcd11690f6d2ad84/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7bf84cc6-557b-4a5e-9688-f355a2432ea0.c0005c834c23b3316dcc6d7d8b2418b1153ebe82.de-de.xlf"
$de05eaXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/098681995e41771b61ad938777136cb16c3ad337/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/05ea4a7e-cab7-444d-9b45-abe39f53d2f7.544af659032ec6d29b8d61d45c87655c2dbfde8c.de-de.xlf"
$de7bf8XlfName = "7bf84cc6-557b-4a5e-9688-f355a2432ea0.c0005c834c23b3316dcc6d7d8b2418b1153ebe82.de-de.xlf"
$de05eaXlfName = "05ea4a7e-cab7-444d-9b45-abe39f53d2f7.544af659032ec6d29b8d61d45c87655c2dbfde8c.de-de.xlf"

$wsDe.Range("A1").Hyperlinks.Delete()

# Row 2: 7bf84cc6 item - handed back, fully populated (incl. new columns).
$wsDe.Range("A2").Value = $md7bf8Name
$wsDe.Range("B2").Value = $handedBack
$wsDe.Range("C2").Value = $de7bf8XlfName
$wsDe.Range("D2").Value = "2016-03-11 00:50:42"
$wsDe.Range("E2").Value = $md7bf8Name
$wsDe.Range("F2").Value = $de7bf8XlfName
$wsDe.Range("G2").Value = "2016-03-11 00:51:52"
$wsDe.Range("H2").Value = "Include"

# Row 3: 05ea4a7e item - still ready for handoff.
$wsDe.Range("A3").Value = $md05eaName
$wsDe.Range("B3").Value = $readyOff
$wsDe.Range("C3").Value = $de05eaXlfName
$wsDe.Range("D3").Value = "2016-03-11 00:50:12"
$wsDe.Range("G3").Value = "0001-01-01 00:00:00"
$wsDe.Range("H3").Value = "Include"

# Row 4: .localization-config - unchanged.
$wsDe.Range("A4").Value = $configName
$wsDe.Range("B4").Value = $notLocalized
$wsDe.Range("D4").Value = "0001-01-01 00:00:00"
$wsDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDe.Range("H4").Value = "Ignored"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $md7bf8Url, "", "", $md7bf8Name)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $de7bf8XlfUrl, "", "", $de7bf8XlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $md7bf8Url, "", "", $md7bf8Name)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $de7bf8XlfUrl, "", "", $de7bf8XlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $md05eaUrl, "", "", $md05eaName)
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), $de05eaXlfUrl, "", "", $de05eaXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $configUrl, "", "", $configName)

foreach ($addr in @("A2", "C2", "E2", "F2", "A3", "C3", "A4")) {
    $wsDe.Range($addr).Style = "HyperLink"
}
$wsDe.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
